$d = $word.ActiveDocument

# --- 1. First paragraph: append two trailing spaces to the existing
#        sentence, then append a new, red-colored parenthetical note.
#        The red note is typed/colored in three chunks so the resulting
#        OOXML run-split matches the authored edit exactly.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs.Item(1)
$endRange = $p1.Range
$endRange.SetRange($endRange.End - 1, $endRange.End - 1)

$enDash = [char]0x2013
$chunk1 = "(This is a change " + $enDash + " Ve"
$chunk2 = "rsion for main branch"
$chunk3 = ")"
$fullNote = $chunk1 + $chunk2 + $chunk3

$endRange.InsertAfter($fullNote)

$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Find.Execute($chunk1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$r2 = $p1.Range
$r2.Find.Execute($chunk2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$r3 = $p1.Range
$r3.Find.Execute($chunk3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Font.Color = 255

# --- 2. Remove the trailing "...ank God almighty, we are free at last."
#        paragraph entirely (it was the only user of the NormalWeb
#        paragraph style in the document body).
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Delete()
